# Apply updated crypto price/volume data to worksheet cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.269.65"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.593.81"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "1.816.76"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "1.586.91"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "26.255.19"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -2.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "1.420.11"
$ws.Range("E33").Value = "  +6.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -6.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.765"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "1.729.71"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("E51").Value = "  -0.04%  "
